$wb = $excel.ActiveWorkbook

# Rename the existing "Personas" sheet to "futbol varones"
$ws1 = $wb.Worksheets.Item("Personas")
$ws1.Name = "futbol varones"

# Add a new worksheet after the first one, named "futbol damas"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "futbol damas"

# Populate header row on the new sheet, matching the first sheet's headers
$headers = @("nombre", "apellido", "email", "rut", "phone_number", "emergency_phone")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}
